$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "who?"
$ws.Range("G2").Value = "me"
$ws.Range("G3").Value = "me"
$ws.Range("G4").Value = "me"
$ws.Range("G5").Value = "me"
$ws.Range("G6").Value = "me"
$ws.Range("G7").Value = "me"
$ws.Range("G8").Value = "me"
$ws.Range("G9").Value = "me"
$ws.Range("G10").Value = "me"

$ws.Range("I4").Select()
